$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-sort the measurement table (A2:F7) in ascending order by column A.
# The data had previously been sorted descending by column A; sorting it
# ascending reorders all the rows and records the applied sort state,
# matching Excel's own Data > Sort behaviour.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A7"))
$ws.Sort.SetRange($ws.Range("A2:F7"))
$ws.Sort.Header = -4142
$ws.Sort.Apply()

# Reflect the sorted table as the active selection.
$ws.Range("A2:E7").Select() | Out-Null
